# HomeLyf data sheet: add "QDTAS1" worksheet with user end-to-end test data,
# and update Sheet1's selection (no more active-cell pin at B2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: drop the old tabSelected/active-cell pin, select header row instead
$ws1.Range("A1:E1").Select()

# --- Add the new worksheet right after Sheet1 and name it
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "QDTAS1"

# --- Header row
$ws2.Range("A1").Value = "userName"
$ws2.Range("B1").Value = "email"
$ws2.Range("C1").Value = "password"
$ws2.Range("D1").Value = "firstName"
$ws2.Range("E1").Value = "middleName"
$ws2.Range("F1").Value = "lastName"
$ws2.Range("G1").Value = "gender"
$ws2.Range("H1").Value = "depId"
$ws2.Range("I1").Value = "role"
$ws2.Range("J1").Value = "phoneNubmer"
$ws2.Range("K1").Value = "designation"
$ws2.Range("L1").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("L1").Value = "birthDate"
$ws2.Range("M1").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("M1").Value = "joinDate"
$ws2.Range("N1").Value = "address"

# --- Row 2 (Onkar Kure a.k.a. "onkark")
$ws2.Range("A2").Value = "onkark"
$ws2.Range("B2").Value = "kure.onkar41@gmail.com"
$ws2.Range("C2").Value = "Onkar@41"
$ws2.Range("D2").Value = "omkar"
$ws2.Range("E2").Value = "shiv"
$ws2.Range("F2").Value = "kure"
$ws2.Range("G2").Value = "male"
$ws2.Range("H2").Value = 2
$ws2.Range("I2").Value = "QA"
$ws2.Range("J2").Value = 1234567890
$ws2.Range("K2").Value = "Automation.Engineer"
$ws2.Range("L2").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("L2").Value = "2024-02-12"
$ws2.Range("M2").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("M2").Value = "2024-02-12"
$ws2.Range("N2").Value = "Pune"

# --- Row 3 (Onkar Singal)
$ws2.Range("A3").Value = "Onkar"
$ws2.Range("B3").Value = "mayivec121@irnini.com"
$ws2.Range("C3").Value = "Bebo@123"
$ws2.Range("D3").Value = "om"
$ws2.Range("E3").Value = "shiv"
$ws2.Range("F3").Value = "singal"
$ws2.Range("G3").Value = "Male"
$ws2.Range("H3").Value = 2
$ws2.Range("I3").Value = "QA"
$ws2.Range("J3").Value = 1234567892
$ws2.Range("K3").Value = "ManualTester"
$ws2.Range("L3").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("L3").Value = "1998-07-18"
$ws2.Range("M3").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("M3").Value = "2018-08-18"
$ws2.Range("N3").Value = "Pune"

# --- Column widths (best-effort visual fit matching the new data)
$ws2.Columns.Item(2).ColumnWidth = 22.33203125
$ws2.Columns.Item(3).ColumnWidth = 10
$ws2.Columns.Item(10).ColumnWidth = 12.6640625
$ws2.Columns.Item(11).ColumnWidth = 18.44140625
$ws2.Columns.Item(12).ColumnWidth = 10.5546875
$ws2.Columns.Item(13).ColumnWidth = 10.33203125

# --- Hyperlinks on the password / email cells (mailto:, mirroring Sheet1's pattern)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:Bebo@123")
$ws2.Range("C3").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:Onkar@41")
$ws2.Range("C2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:kure.onkar41@gmail.com")
$ws2.Range("B2").Style = "Hyperlink"

# --- Final selection / active cell on the new sheet
$ws2.Range("L4").Select()

Write-Host "QDTAS1 sheet created with end-to-end user data"
